$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1 "time_taken", matching formatting of the existing
# header cells (bold + bordered + centered) by copying E1's format.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("F1").Value = "time_taken"

# New time_taken values for each data row
$ws.Cells.Item(2, 6).Value = "2021-10-05 10:51:38.938661"
$ws.Cells.Item(3, 6).Value = "2021-10-05 10:51:38.938671"
$ws.Cells.Item(4, 6).Value = "2021-10-05 10:51:38.938674"
